# Fruta / hortaliza, semanal
# Insert the latest weekly price record for
# "Femacal de La Calera - Zanahoria" just above the existing entry that is
# currently in row 170 (date 2020-11-26 / serial 44161), pushing every
# following row of this data block down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 170:188 down to 171:189, leaving a blank row 170 to fill in.
$ws.Rows("170:170").Insert()

# Populate the new weekly record in row 170.
$ws.Range("A170").Value = 3
$ws.Range("B170").Value = "Femacal de La Calera"
$ws.Range("C170").Value = "Coquimbo"
$ws.Range("D170").Value2 = 44449
$ws.Range("E170").Value = 5
$ws.Range("F170").Value = 100114013
$ws.Range("G170").Value = "Zanahoria"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 370
$ws.Range("K170").Value = 5000
$ws.Range("L170").Value = 5500
$ws.Range("M170").Value = 5257
$ws.Range("N170").Value = "$/saco 20 kilos"
$ws.Range("O170").Value = "Provincia de Quillota"
$ws.Range("P170").Value = 263
$ws.Range("Q170").Value = 20
$ws.Range("R170").Value = "Hortaliza"
